# Add " for SENG8040" right after the existing "This is lab1 " text.
#
# In the target document the added text lands in its own <w:r> rather
# than being folded into the run that already holds "This is lab1 ",
# even though both runs end up with identical run properties
# (just <w:lang w:val="en-US"/>). We reproduce that by locating the
# insertion point with Find, inserting the new text there, and then
# "touching" a direct-formatting property (Bold on, then back off) on
# the freshly inserted range so it is written out as a run of its own
# instead of silently re-merging with its identically-formatted
# neighbour.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("This is lab1 ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)                    # wdCollapseEnd: move to right after "This is lab1 "
    $rng.InsertAfter(" for SENG8040")   # type the new text

    $rng.Bold = 1
    $rng.Bold = 0                       # force a distinct run without leaving formatting behind
}
